$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update existing value in G5 (model size for MHLA row)
$ws.Range("G5").Value = 16939008

# Add new row 9: "Flash MHLA" results, mirroring the structure of row 8 (MQA / RoPE / latent dim = 32)
$ws.Range("A9").Value = "Flash MHLA"
$ws.Range("B9").Value = "RoPE"
$ws.Range("C9").Value = 8
$ws.Range("D9").Value = 256
$ws.Range("E9").Value = 6
$ws.Range("F9").Value = "latent dim = 32"
$ws.Range("H9").Value = 8
$ws.Range("I9").Value = 1024
$ws.Range("J9").Value = 20000000000000
$ws.Range("J9").NumberFormat = $ws.Range("J8").NumberFormat
$ws.Range("K9").Value = 2500
$ws.Range("N9").Value = 5000

# Column A now needs a best-fit width to accommodate "Flash MHLA"
$ws.Columns.Item(1).AutoFit() | Out-Null
$ws.Columns.Item(1).ColumnWidth = 9.6

# Update the active selection on the sheet
$ws.Range("G15").Select()
